{"js": "// Update the two-digit-division answer table: each cell's text is replaced\n// with a new division fact, per the commit's diff. Old strings are unique\n// across the document, so a plain search+replace for each pair is safe.\nconst replacements = [\n  [\"42\u00f76=7, 0\", \"79\u00f78=9, 7\"],\n  [\"47\u00f75=9, 2\", \"46\u00f75=9, 1\"],\n  [\"79\u00f75=15, 4\", \"77\u00f72=38, 1\"],\n  [\"46\u00f76=7, 4\", \"93\u00f76=15, 3\"],\n  [\"37\u00f72=18, 1\", \"70\u00f72=35, 0\"],\n  [\"14\u00f79=1, 5\", \"21\u00f76=3, 3\"],\n  [\"44\u00f79=4, 8\", \"83\u00f75=16, 3\"],\n  [\"44\u00f74=11, 0\", \"34\u00f78=4, 2\"],\n  [\"42\u00f73=14, 0\", \"11\u00f78=1, 3\"],\n  [\"15\u00f76=2, 3\", \"26\u00f74=6, 2\"],\n  [\"86\u00f74=21, 2\", \"17\u00f77=2, 3\"],\n  [\"33\u00f73=11, 0\", \"30\u00f74=7, 2\"],\n  [\"56\u00f78=7, 0\", \"21\u00f73=7, 0\"],\n  [\"83\u00f72=41, 1\", \"10\u00f79=1, 1\"],\n  [\"93\u00f77=13, 2\", \"78\u00f77=11, 1\"],\n  [\"42\u00f79=4, 6\", \"54\u00f76=9, 0\"],\n  [\"12\u00f79=1, 3\", \"36\u00f79=4, 0\"],\n  [\"64\u00f76=10, 4\", \"92\u00f72=46, 0\"],\n  [\"56\u00f79=6, 2\", \"44\u00f75=8, 4\"],\n  [\"42\u00f75=8, 2\", \"42\u00f72=21, 0\"],\n  [\"18\u00f74=4, 2\", \"70\u00f76=11, 4\"],\n  [\"20\u00f77=2, 6\", \"32\u00f74=8, 0\"],\n  [\"75\u00f79=8, 3\", \"38\u00f75=7, 3\"],\n  [\"17\u00f75=3, 2\", \"36\u00f77=5, 1\"],\n  [\"37\u00f76=6, 1\", \"90\u00f76=15, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit-division answer table: each cell's text is replaced\n# with a new division fact, per the commit's diff. Old strings are unique\n# across the document, so Find/Replace for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"42\u00f76=7, 0\", \"79\u00f78=9, 7\"),\n  @(\"47\u00f75=9, 2\", \"46\u00f75=9, 1\"),\n  @(\"79\u00f75=15, 4\", \"77\u00f72=38, 1\"),\n  @(\"46\u00f76=7, 4\", \"93\u00f76=15, 3\"),\n  @(\"37\u00f72=18, 1\", \"70\u00f72=35, 0\"),\n  @(\"14\u00f79=1, 5\", \"21\u00f76=3, 3\"),\n  @(\"44\u00f79=4, 8\", \"83\u00f75=16, 3\"),\n  @(\"44\u00f74=11, 0\", \"34\u00f78=4, 2\"),\n  @(\"42\u00f73=14, 0\", \"11\u00f78=1, 3\"),\n  @(\"15\u00f76=2, 3\", \"26\u00f74=6, 2\"),\n  @(\"86\u00f74=21, 2\", \"17\u00f77=2, 3\"),\n  @(\"33\u00f73=11, 0\", \"30\u00f74=7, 2\"),\n  @(\"56\u00f78=7, 0\", \"21\u00f73=7, 0\"),\n  @(\"83\u00f72=41, 1\", \"10\u00f79=1, 1\"),\n  @(\"93\u00f77=13, 2\", \"78\u00f77=11, 1\"),\n  @(\"42\u00f79=4, 6\", \"54\u00f76=9, 0\"),\n  @(\"12\u00f79=1, 3\", \"36\u00f79=4, 0\"),\n  @(\"64\u00f76=10, 4\", \"92\u00f72=46, 0\"),\n  @(\"56\u00f79=6, 2\", \"44\u00f75=8, 4\"),\n  @(\"42\u00f75=8, 2\", \"42\u00f72=21, 0\"),\n  @(\"18\u00f74=4, 2\", \"70\u00f76=11, 4\"),\n  @(\"20\u00f77=2, 6\", \"32\u00f74=8, 0\"),\n  @(\"75\u00f79=8, 3\", \"38\u00f75=7, 3\"),\n  @(\"17\u00f75=3, 2\", \"36\u00f77=5, 1\"),\n  @(\"37\u00f76=6, 1\", \"90\u00f76=15, 0\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $old\n  $find.Replacement.Text = $new\n  $find.Execute(\n    $old,\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    0,\n    $false,\n    $new,\n    2\n  )\n}\n\nWrite-Host \"Replaced $($pairs.Count) division facts\"\n"}
